# Proofreading edits to the non-technical MCM summary.
#
# Paragraph 2 (the "While manufacturers..." paragraph) is rewritten to open
# with the new "stringently regulated industries" framing, and the _GoBack
# bookmark (previously sitting alone in its own paragraph near the end of the
# document) is relocated to its new "last edit" position inside this
# paragraph.
#
# Paragraph 3 (the "However, a mathematical model..." paragraph) is reworked
# with a couple of still-unfinished sentences (flagged by Word's grammar
# checker via <w:proofErr> gramStart/gramEnd markers) plus several smaller
# wording tweaks further down.
#
# Both paragraphs are rebuilt with InsertXML so the resulting run layout and
# bookmark/proofErr placement match exactly.

$d = $word.ActiveDocument

# --- Remove the old stand-alone _GoBack bookmark paragraph near the end ---
# (do this first, before a new _GoBack is (re)inserted elsewhere below, so
# the lookup unambiguously finds the original bookmark's home paragraph).
# It now becomes a plain empty paragraph (the bookmark itself moves into
# paragraph 2 below, to mark the new "last edit" position).
$emptyXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.WordOpenXML -like "*_GoBack*") {
        $para.Range.InsertXML($emptyXml)
        break
    }
}

# --- Paragraph 2: "While manufacturers..." -> "Aircraft manufacturing..." ---
$p1xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:r><w:t>Aircraft manufacturing and air traffic control are among the most stringently regulated industries in the world</w:t></w:r>
<w:r><w:t>.</w:t></w:r>
<w:r><w:t xml:space="preserve"> Since this is the case, we take the issue of losing aircraft over oceans extremely seriously. </w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
<w:r><w:t>Locating a plane of any size in any of the vast oceans on our planet is very challenging.</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t>E</w:t></w:r>
<w:r><w:t>ven</w:t></w:r>
<w:r><w:t xml:space="preserve"> when debris of a downed aircraft has been located,</w:t></w:r>
<w:r><w:t xml:space="preserve"> finding </w:t></w:r>
<w:r><w:t>where the rest of that airplane resides remains a daunting task.</w:t></w:r>
</w:p>
'@
$d.Paragraphs.Item(2).Range.InsertXML($p1xml)

# --- Paragraph 3: "However, a mathematical model..." rewording ---
$p2xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:r><w:t xml:space="preserve">In order to </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:t xml:space="preserve">address </w:t></w:r>
<w:r><w:t>.</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t xml:space="preserve"> By incorporating probabili</w:t></w:r>
<w:r><w:t>ty theory and well</w:t></w:r>
<w:r><w:t xml:space="preserve"> researched mathematical theories, t</w:t></w:r>
<w:r><w:t>his model optimizes</w:t></w:r>
<w:r><w:t xml:space="preserve"> the</w:t></w:r>
<w:r><w:t xml:space="preserve"> probability of locating the aircraft w</w:t></w:r>
<w:r><w:t>hile also minimizing the cost of</w:t></w:r>
<w:r><w:t xml:space="preserve"> search</w:t></w:r>
<w:r><w:t>ing</w:t></w:r>
<w:r><w:t xml:space="preserve">. </w:t></w:r>
<w:r><w:t>Fortunately, this model</w:t></w:r>
<w:r><w:t xml:space="preserve"> has been developed to be </w:t></w:r>
<w:r><w:t xml:space="preserve">adapted </w:t></w:r>
<w:r><w:t>for different needs</w:t></w:r>
<w:r><w:t xml:space="preserve"> and extendable past its </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:t xml:space="preserve">current </w:t></w:r>
<w:r><w:t>.</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t xml:space="preserve"> It is</w:t></w:r>
<w:r><w:t xml:space="preserve"> possible to adjust multiple</w:t></w:r>
<w:r><w:t xml:space="preserve"> parameters that </w:t></w:r>
<w:r><w:t>technologies being used</w:t></w:r>
<w:r><w:t>,</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t xml:space="preserve">the </w:t></w:r>
<w:r><w:t>number of search aircraft</w:t></w:r>
<w:r><w:t xml:space="preserve"> being deployed, and the breadth of area being searched</w:t></w:r>
<w:r><w:t>, making it applicable in a</w:t></w:r>
<w:r><w:t xml:space="preserve"> multitude of scenarios.</w:t></w:r>
</w:p>
'@
$d.Paragraphs.Item(3).Range.InsertXML($p2xml)

Write-Host "Done"
